# Reorder the invoice report columns (C..K) into their new layout and
# resize the affected columns to match.
#
# Old layout (cols C..K): GSTIN, Invoice No., Date, Taxable Amount, Total Tax,
#                          CGST, SGST, IGST, HSN Codes
# New layout (cols C..K): Date, GSTIN, Invoice No., HSN Codes, CGST, SGST,
#                          IGST, Total Tax, Taxable Amount
#
# Values/styles are moved with native Range.Copy(Destination) calls (instead
# of re-typing .Value) so that text cells that merely look like numbers or
# dates (e.g. "2024-09-28", "173.91") are relocated verbatim, without Excel
# reinterpreting them as real dates/numbers, and so each cell keeps the
# formatting that belongs to its content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Stage the current C1:K2 block far away from the live data so source
#     and destination ranges never overlap while we shuffle cells around. ---
$staging = $ws.Range("C100")
$ws.Range("C1:K2").Copy($staging)

# Row 1 (headers) staged at row 100, row 2 (data) staged at row 101.
# Old column order: C=GSTIN D=Invoice No. E=Date F=Taxable Amount G=Total Tax
#                    H=CGST I=SGST J=IGST K=HSN Codes
$s_gstin_h   = $ws.Range("C100")
$s_invoice_h = $ws.Range("D100")
$s_date_h    = $ws.Range("E100")
$s_taxable_h = $ws.Range("F100")
$s_total_h   = $ws.Range("G100")
$s_cgst_h    = $ws.Range("H100")
$s_sgst_h    = $ws.Range("I100")
$s_igst_h    = $ws.Range("J100")
$s_hsn_h     = $ws.Range("K100")

$s_gstin_d   = $ws.Range("C101")
$s_invoice_d = $ws.Range("D101")
$s_date_d    = $ws.Range("E101")
$s_taxable_d = $ws.Range("F101")
$s_total_d   = $ws.Range("G101")
$s_cgst_d    = $ws.Range("H101")
$s_sgst_d    = $ws.Range("I101")
$s_igst_d    = $ws.Range("J101")
$s_hsn_d     = $ws.Range("K101")

# --- Copy back into the live range in the new column order ---
$s_date_h.Copy($ws.Range("C1"))
$s_gstin_h.Copy($ws.Range("D1"))
$s_invoice_h.Copy($ws.Range("E1"))
$s_hsn_h.Copy($ws.Range("F1"))
$s_cgst_h.Copy($ws.Range("G1"))
$s_sgst_h.Copy($ws.Range("H1"))
$s_igst_h.Copy($ws.Range("I1"))
$s_total_h.Copy($ws.Range("J1"))
$s_taxable_h.Copy($ws.Range("K1"))

$s_date_d.Copy($ws.Range("C2"))
$s_gstin_d.Copy($ws.Range("D2"))
$s_invoice_d.Copy($ws.Range("E2"))
$s_hsn_d.Copy($ws.Range("F2"))
$s_cgst_d.Copy($ws.Range("G2"))
$s_sgst_d.Copy($ws.Range("H2"))
$s_igst_d.Copy($ws.Range("I2"))
$s_total_d.Copy($ws.Range("J2"))
$s_taxable_d.Copy($ws.Range("K2"))

# --- Remove the temporary staging area entirely (content + formatting) ---
$ws.Range("C100:K101").Clear()

# --- Resize the columns to match the new layout ---
$ws.Columns.Item(3).ColumnWidth = 14.17   # C: 15
$ws.Columns.Item(4).ColumnWidth = 17.17   # D: 18
$ws.Columns.Item(5).ColumnWidth = 19.17   # E: 20
$ws.Columns.Item(6).ColumnWidth = 39.17   # F: 40
$ws.Columns.Item(11).ColumnWidth = 14.17  # K: 15
